# Update division problems to new values per commit diff.
$d = $word.ActiveDocument

$replacements = @(
    @("61÷3=", "76÷4="),
    @("25÷2=", "38÷9="),
    @("23÷4=", "19÷7="),
    @("65÷4=", "40÷8="),
    @("94÷8=", "19÷6="),
    @("39÷7=", "53÷9="),
    @("78÷8=", "50÷7="),
    @("56÷9=", "81÷4="),
    @("68÷5=", "56÷7="),
    @("53÷6=", "41÷4="),
    @("79÷4=", "75÷7="),
    @("29÷7=", "88÷6="),
    @("59÷8=", "20÷4="),
    @("72÷7=", "80÷6="),
    @("24÷2=", "89÷4="),
    @("27÷4=", "89÷8="),
    @("47÷7=", "86÷2="),
    @("65÷7=", "14÷2="),
    @("79÷8=", "21÷9="),
    @("96÷4=", "11÷5="),
    @("81÷9=", "72÷4="),
    @("14÷3=", "83÷7="),
    @("96÷8=", "26÷5="),
    @("82÷9=", "38÷2=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
